$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 26 data: Leetcode Question No. 160, "Intersection of Two Linked Lists", "Java,Python"
$ws.Cells.Item(26, 1).Value = 160
$ws.Cells.Item(26, 2).Value = "Intersection of Two Linked Lists"
$ws.Cells.Item(26, 3).Value = "Java,Python"

# Match style of the A column cells (left/top aligned) as used by neighboring rows (e.g. A25)
$ws.Cells.Item(26, 1).HorizontalAlignment = -4131  # xlLeft
$ws.Cells.Item(26, 1).VerticalAlignment = -4160    # xlTop

# Update the view: scroll so row 7 is the top-left visible row, and select D26
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("D26").Select()
